$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits inside the
#    "PhD Biomedical engineering" paragraph (between the "B" run and
#    the "iomedical engineering" run).
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2) Insert a new "File Versioning Control: Git" paragraph right after
#    the "Research Methods: ..." paragraph and before the "Markup and
#    Typesetting Languages: ..." paragraph, then re-anchor the
#    _GoBack bookmark (collapsed) at the very end of the new paragraph.
# ---------------------------------------------------------------------
$researchMethods = $d.Paragraphs.Item(39)
$researchMethods.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(40)
$insertRange = $newPara.Range.Duplicate
$insertRange.Collapse(1)

# A trailing placeholder "|" marks where the paragraph's real content
# ends (before its own paragraph mark) so that the bookmark can be
# dropped in a safely mid-paragraph position; the placeholder is
# removed again right afterwards.
$insertRange.InsertBefore("File Versioning Control: Git|")

$paraStart = $newPara.Range.Start
$headingLen = "File Versioning Control".Length
$headingRange = $d.Range($paraStart, $paraStart + $headingLen)
$headingRange.Font.Bold = 1
$headingRange.Font.Underline = 1

$bodyLenNoPlaceholder = "File Versioning Control: Git".Length
$placeholderPos = $paraStart + $bodyLenNoPlaceholder

$bmRange = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()
